$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 39427
$ws.Range("D2").Value = 56990721
$ws.Range("C3").Value = 94365
$ws.Range("D3").Value = 138299236
$ws.Range("C4").Value = 32179
$ws.Range("D4").Value = 47642976
$ws.Range("C5").Value = 9071
$ws.Range("D5").Value = 13481753
$ws.Range("C6").Value = 2134
$ws.Range("D6").Value = 3172971
$ws.Range("C7").Value = 181
$ws.Range("D7").Value = 266593
$ws.Range("C12").Value = 42803
$ws.Range("D12").Value = 58009257
$ws.Range("C13").Value = 10052
$ws.Range("D13").Value = 14527718
$ws.Range("C14").Value = 26794
$ws.Range("D14").Value = 39274530
$ws.Range("C15").Value = 8555
$ws.Range("D15").Value = 12696643
$ws.Range("C16").Value = 2234
$ws.Range("D16").Value = 3318849
$ws.Range("C17").Value = 439
$ws.Range("D17").Value = 647623
$ws.Range("C18").Value = 38
$ws.Range("D18").Value = 57000
$ws.Range("C20").Value = 10531
$ws.Range("D20").Value = 13911985
$ws.Range("C21").Value = 13881
$ws.Range("D21").Value = 20022026
$ws.Range("C22").Value = 32600
$ws.Range("D22").Value = 47821836
$ws.Range("C23").Value = 10526
$ws.Range("D23").Value = 15643646
$ws.Range("C24").Value = 2738
$ws.Range("D24").Value = 4071674
$ws.Range("C25").Value = 556
$ws.Range("D25").Value = 828092
$ws.Range("C26").Value = 39
$ws.Range("D26").Value = 57953
$ws.Range("C27").Value = 12054
$ws.Range("D27").Value = 16072402
$ws.Range("C28").Value = 7998
$ws.Range("D28").Value = 11573884
$ws.Range("C29").Value = 23284
$ws.Range("D29").Value = 34173049
$ws.Range("C30").Value = 8024
$ws.Range("D30").Value = 11930241
$ws.Range("C31").Value = 2038
$ws.Range("D31").Value = 3041251
$ws.Range("C32").Value = 382
$ws.Range("D32").Value = 570415
$ws.Range("C34").Value = 8597
$ws.Range("D34").Value = 11355038
$ws.Range("C35").Value = 3411
$ws.Range("D35").Value = 4925201
$ws.Range("C36").Value = 8173
$ws.Range("D36").Value = 11939222
$ws.Range("C37").Value = 3263
$ws.Range("D37").Value = 4837461
$ws.Range("C38").Value = 845
$ws.Range("D38").Value = 1258555
$ws.Range("C41").Value = 2574
$ws.Range("D41").Value = 3478381
$ws.Range("C42").Value = 17968
$ws.Range("D42").Value = 25965762
$ws.Range("C43").Value = 52808
$ws.Range("D43").Value = 77387689
$ws.Range("C44").Value = 19467
$ws.Range("D44").Value = 28904864
$ws.Range("C45").Value = 5815
$ws.Range("D45").Value = 8654435
$ws.Range("C46").Value = 1289
$ws.Range("D46").Value = 1923144
$ws.Range("C47").Value = 70
$ws.Range("D47").Value = 103068
$ws.Range("C50").Value = 17335
$ws.Range("D50").Value = 23011303
$ws.Range("C51").Value = 2196
$ws.Range("D51").Value = 3188096
$ws.Range("C52").Value = 7399
$ws.Range("D52").Value = 10872534
$ws.Range("C53").Value = 2471
$ws.Range("D53").Value = 3688644
$ws.Range("C54").Value = 779
$ws.Range("D54").Value = 1163415
$ws.Range("C55").Value = 204
$ws.Range("D55").Value = 302226
$ws.Range("C57").Value = 7487
$ws.Range("D57").Value = 10296246
$ws.Range("C58").Value = 1294
$ws.Range("D58").Value = 2323507
$ws.Range("C59").Value = 3160
$ws.Range("D59").Value = 5665955
$ws.Range("C60").Value = 1258
$ws.Range("D60").Value = 2275566
$ws.Range("C61").Value = 419
$ws.Range("D61").Value = 747083
$ws.Range("C62").Value = 142
$ws.Range("D62").Value = 265100
$ws.Range("C64").Value = 1933
$ws.Range("D64").Value = 3205583
$ws.Range("C65").Value = 16076
$ws.Range("D65").Value = 23218312
$ws.Range("C66").Value = 46326
$ws.Range("D66").Value = 67755484
$ws.Range("C67").Value = 16172
$ws.Range("D67").Value = 24027361
$ws.Range("C68").Value = 4730
$ws.Range("D68").Value = 7044788
$ws.Range("C69").Value = 992
$ws.Range("D69").Value = 1474802
$ws.Range("C73").Value = 15586
$ws.Range("D73").Value = 20510197
$ws.Range("C74").Value = 56359
$ws.Range("D74").Value = 81971960
$ws.Range("C75").Value = 156406
$ws.Range("D75").Value = 230344491
$ws.Range("C76").Value = 67144
$ws.Range("D76").Value = 100025650
$ws.Range("C77").Value = 21573
$ws.Range("D77").Value = 32232013
$ws.Range("C78").Value = 5202
$ws.Range("D78").Value = 7771348
$ws.Range("C79").Value = 305
$ws.Range("D79").Value = 452670
$ws.Range("C80").Value = 27
$ws.Range("D80").Value = 39405
$ws.Range("C85").Value = 55128
$ws.Range("D85").Value = 74784341
$ws.Range("C86").Value = 4849
$ws.Range("D86").Value = 7027173
$ws.Range("C87").Value = 12026
$ws.Range("D87").Value = 17665491
$ws.Range("C88").Value = 3998
$ws.Range("D88").Value = 5956958
$ws.Range("C89").Value = 1384
$ws.Range("D89").Value = 2067111
$ws.Range("C90").Value = 304
$ws.Range("D90").Value = 453512
$ws.Range("C93").Value = 5638
$ws.Range("D93").Value = 7573062
$ws.Range("C94").Value = 1690
$ws.Range("D94").Value = 2436567
$ws.Range("C95").Value = 5454
$ws.Range("D95").Value = 8036497
$ws.Range("C96").Value = 2012
$ws.Range("D96").Value = 2994808
$ws.Range("C97").Value = 725
$ws.Range("D97").Value = 1086460
$ws.Range("C98").Value = 202
$ws.Range("D98").Value = 305113
$ws.Range("C99").Value = 21
$ws.Range("D99").Value = 31500
$ws.Range("C101").Value = 3750
$ws.Range("D101").Value = 4971238
$ws.Range("C102").Value = 781
$ws.Range("D102").Value = 1380245
$ws.Range("C103").Value = 484
$ws.Range("D103").Value = 885004
$ws.Range("C104").Value = 185
$ws.Range("D104").Value = 339789
$ws.Range("C105").Value = 62
$ws.Range("D105").Value = 112500
$ws.Range("C107").Value = 11277
$ws.Range("D107").Value = 16353792
$ws.Range("C108").Value = 30099
$ws.Range("D108").Value = 44197846
$ws.Range("C109").Value = 10082
$ws.Range("D109").Value = 14988726
$ws.Range("C110").Value = 2776
$ws.Range("D110").Value = 4139080
$ws.Range("C111").Value = 513
$ws.Range("D111").Value = 764546
$ws.Range("C112").Value = 60
$ws.Range("D112").Value = 90000
$ws.Range("C114").Value = 10116
$ws.Range("D114").Value = 13338540
$ws.Range("C115").Value = 31713
$ws.Range("D115").Value = 45710108
$ws.Range("C116").Value = 68291
$ws.Range("D116").Value = 99914945
$ws.Range("C117").Value = 21937
$ws.Range("D117").Value = 32593375
$ws.Range("C118").Value = 6240
$ws.Range("D118").Value = 9295299
$ws.Range("C119").Value = 1184
$ws.Range("D119").Value = 1768889
$ws.Range("C120").Value = 91
$ws.Range("D120").Value = 132895
$ws.Range("C124").Value = 26627
$ws.Range("D124").Value = 35524574
$ws.Range("C125").Value = 37625
$ws.Range("D125").Value = 54277488
$ws.Range("C126").Value = 79598
$ws.Range("D126").Value = 116363478
$ws.Range("C127").Value = 24595
$ws.Range("D127").Value = 36502859
$ws.Range("C128").Value = 6598
$ws.Range("D128").Value = 9805623
$ws.Range("C129").Value = 1324
$ws.Range("D129").Value = 1969311
$ws.Range("C130").Value = 70
$ws.Range("D130").Value = 103228
$ws.Range("C133").Value = 32876
$ws.Range("D133").Value = 43622567
$ws.Range("C134").Value = 13835
$ws.Range("D134").Value = 20027193
$ws.Range("C135").Value = 33398
$ws.Range("D135").Value = 49042832
$ws.Range("C136").Value = 11791
$ws.Range("D136").Value = 17519816
$ws.Range("C137").Value = 3085
$ws.Range("D137").Value = 4598741
$ws.Range("C138").Value = 523
$ws.Range("D138").Value = 778490
$ws.Range("C139").Value = 40
$ws.Range("D139").Value = 58825
$ws.Range("C141").Value = 11161
$ws.Range("D141").Value = 14869749
$ws.Range("C142").Value = 36780
$ws.Range("D142").Value = 53119939
$ws.Range("C143").Value = 84670
$ws.Range("D143").Value = 124037636
$ws.Range("C144").Value = 25177
$ws.Range("D144").Value = 37403613
$ws.Range("C145").Value = 6624
$ws.Range("D145").Value = 9884496
$ws.Range("C146").Value = 1508
$ws.Range("D146").Value = 2244230
$ws.Range("C149").Value = 30279
$ws.Range("D149").Value = 40800310
